# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Monsters")
$ws2 = $wb.Worksheets.Item("Monsters Skills")

# --- Sheet 'Monsters': append new monster rows (43-57) ---
$ws1.Cells.Item(43,1).Value = 'Adult Red Dragon'
$ws1.Cells.Item(43,2).Value = 4500
$ws1.Cells.Item(43,3).Value = 6000
$ws1.Cells.Item(43,4).Value = 500
$ws1.Cells.Item(43,5).Value = 6300
$ws1.Cells.Item(43,6).Value = 6400
$ws1.Cells.Item(43,7).Value = 450
$ws1.Cells.Item(43,8).Value = 110
$ws1.Cells.Item(43,9).Value = 'int'
$ws1.Cells.Item(43,10).Value = 10
$ws1.Cells.Item(43,11).Value = 0.02
$ws1.Cells.Item(43,12).Value = 1000
$ws1.Cells.Item(43,13).Value = '9000-14000'
$ws1.Cells.Item(43,14).Value = '4500-6000'
$ws1.Cells.Item(43,15).Value = 1
$ws1.Cells.Item(43,18).Value = 'Labyrinth'

$ws1.Cells.Item(44,1).Value = 'Litch Lord of the Labyrinth'
$ws1.Cells.Item(44,2).Value = 8000
$ws1.Cells.Item(44,3).Value = 9000
$ws1.Cells.Item(44,4).Value = 7560
$ws1.Cells.Item(44,5).Value = 9700
$ws1.Cells.Item(44,6).Value = 9100
$ws1.Cells.Item(44,7).Value = 475
$ws1.Cells.Item(44,8).Value = 160
$ws1.Cells.Item(44,9).Value = 'chr'
$ws1.Cells.Item(44,10).Value = 10
$ws1.Cells.Item(44,11).Value = 0.03
$ws1.Cells.Item(44,12).Value = 1200
$ws1.Cells.Item(44,13).Value = '14000-17000'
$ws1.Cells.Item(44,14).Value = '6200-8000'
$ws1.Cells.Item(44,15).Value = 1
$ws1.Cells.Item(44,18).Value = 'Labyrinth'

$ws1.Cells.Item(45,1).Value = 'Whisper of Fate'
$ws1.Cells.Item(45,2).Value = 10000
$ws1.Cells.Item(45,3).Value = 10600
$ws1.Cells.Item(45,4).Value = 11000
$ws1.Cells.Item(45,5).Value = 10870
$ws1.Cells.Item(45,6).Value = 12000
$ws1.Cells.Item(45,7).Value = 500
$ws1.Cells.Item(45,8).Value = 200
$ws1.Cells.Item(45,9).Value = 'int'
$ws1.Cells.Item(45,10).Value = 10
$ws1.Cells.Item(45,11).Value = 0.05
$ws1.Cells.Item(45,12).Value = 1300
$ws1.Cells.Item(45,13).Value = '19000-25000'
$ws1.Cells.Item(45,14).Value = '9000-12000'
$ws1.Cells.Item(45,15).Value = 1
$ws1.Cells.Item(45,18).Value = 'Labyrinth'

$ws1.Cells.Item(46,1).Value = 'Demonic Angel'
$ws1.Cells.Item(46,2).Value = 14900
$ws1.Cells.Item(46,3).Value = 13000
$ws1.Cells.Item(46,4).Value = 14999
$ws1.Cells.Item(46,5).Value = 14670
$ws1.Cells.Item(46,6).Value = 14500
$ws1.Cells.Item(46,7).Value = 560
$ws1.Cells.Item(46,8).Value = 253
$ws1.Cells.Item(46,9).Value = 'dex'
$ws1.Cells.Item(46,10).Value = 10
$ws1.Cells.Item(46,11).Value = 0.03
$ws1.Cells.Item(46,12).Value = 1600
$ws1.Cells.Item(46,13).Value = '25000-30000'
$ws1.Cells.Item(46,14).Value = '12300-14500'
$ws1.Cells.Item(46,15).Value = 1
$ws1.Cells.Item(46,18).Value = 'Labyrinth'

$ws1.Cells.Item(47,1).Value = 'Virgin Priestess Zombie'
$ws1.Cells.Item(47,2).Value = 16000
$ws1.Cells.Item(47,3).Value = 14600
$ws1.Cells.Item(47,4).Value = 15999
$ws1.Cells.Item(47,5).Value = 16000
$ws1.Cells.Item(47,6).Value = 16000
$ws1.Cells.Item(47,7).Value = 620
$ws1.Cells.Item(47,8).Value = 275
$ws1.Cells.Item(47,9).Value = 'str'
$ws1.Cells.Item(47,10).Value = 10
$ws1.Cells.Item(47,11).Value = 0.04
$ws1.Cells.Item(47,12).Value = 1800
$ws1.Cells.Item(47,13).Value = '32000-35000'
$ws1.Cells.Item(47,14).Value = '15300-18900'
$ws1.Cells.Item(47,15).Value = 1
$ws1.Cells.Item(47,18).Value = 'Labyrinth'

$ws1.Cells.Item(48,1).Value = 'Hells Paladin'
$ws1.Cells.Item(48,2).Value = 18750
$ws1.Cells.Item(48,3).Value = 18500
$ws1.Cells.Item(48,4).Value = 19000
$ws1.Cells.Item(48,5).Value = 18900
$ws1.Cells.Item(48,6).Value = 18950
$ws1.Cells.Item(48,7).Value = 660
$ws1.Cells.Item(48,8).Value = 325
$ws1.Cells.Item(48,9).Value = 'dex'
$ws1.Cells.Item(48,10).Value = 10
$ws1.Cells.Item(48,11).Value = 0.01
$ws1.Cells.Item(48,12).Value = 2000
$ws1.Cells.Item(48,13).Value = '35000-41000'
$ws1.Cells.Item(48,14).Value = '19000-22500'
$ws1.Cells.Item(48,15).Value = 1
$ws1.Cells.Item(48,18).Value = 'Labyrinth'

$ws1.Cells.Item(49,1).Value = 'Ancient Gold Dragon'
$ws1.Cells.Item(49,2).Value = 22000
$ws1.Cells.Item(49,3).Value = 22040
$ws1.Cells.Item(49,4).Value = 21000
$ws1.Cells.Item(49,5).Value = 21000
$ws1.Cells.Item(49,6).Value = 2100
$ws1.Cells.Item(49,7).Value = 690
$ws1.Cells.Item(49,8).Value = 375
$ws1.Cells.Item(49,9).Value = 'int'
$ws1.Cells.Item(49,10).Value = 10
$ws1.Cells.Item(49,11).Value = 0.04
$ws1.Cells.Item(49,12).Value = 2200
$ws1.Cells.Item(49,13).Value = '41500-4700'
$ws1.Cells.Item(49,14).Value = '23000-25000'
$ws1.Cells.Item(49,15).Value = 1
$ws1.Cells.Item(49,18).Value = 'Labyrinth'

$ws1.Cells.Item(50,1).Value = 'Broken Hearted Banshee'
$ws1.Cells.Item(50,2).Value = 23000
$ws1.Cells.Item(50,3).Value = 21000
$ws1.Cells.Item(50,4).Value = 22000
$ws1.Cells.Item(50,5).Value = 24000
$ws1.Cells.Item(50,6).Value = 20000
$ws1.Cells.Item(50,7).Value = 715
$ws1.Cells.Item(50,8).Value = 400
$ws1.Cells.Item(50,9).Value = 'chr'
$ws1.Cells.Item(50,10).Value = 10
$ws1.Cells.Item(50,11).Value = 0.03
$ws1.Cells.Item(50,12).Value = 2400
$ws1.Cells.Item(50,13).Value = '47000-50000'
$ws1.Cells.Item(50,14).Value = '25500-27000'
$ws1.Cells.Item(50,15).Value = 1
$ws1.Cells.Item(50,18).Value = 'Labyrinth'

$ws1.Cells.Item(51,1).Value = 'Red Wizard of Krull'
$ws1.Cells.Item(51,2).Value = 24000
$ws1.Cells.Item(51,3).Value = 23400
$ws1.Cells.Item(51,4).Value = 25000
$ws1.Cells.Item(51,5).Value = 24500
$ws1.Cells.Item(51,6).Value = 2500
$ws1.Cells.Item(51,7).Value = 745
$ws1.Cells.Item(51,8).Value = 440
$ws1.Cells.Item(51,9).Value = 'dur'
$ws1.Cells.Item(51,10).Value = 10
$ws1.Cells.Item(51,11).Value = 0.01
$ws1.Cells.Item(51,12).Value = 2600
$ws1.Cells.Item(51,13).Value = '50000-54000'
$ws1.Cells.Item(51,14).Value = '28000-31500'
$ws1.Cells.Item(51,15).Value = 1
$ws1.Cells.Item(51,18).Value = 'Labyrinth'

$ws1.Cells.Item(52,1).Value = 'Devils Lover'
$ws1.Cells.Item(52,2).Value = 25000
$ws1.Cells.Item(52,3).Value = 27000
$ws1.Cells.Item(52,4).Value = 26000
$ws1.Cells.Item(52,5).Value = 25000
$ws1.Cells.Item(52,6).Value = 28600
$ws1.Cells.Item(52,7).Value = 800
$ws1.Cells.Item(52,8).Value = 490
$ws1.Cells.Item(52,9).Value = 'int'
$ws1.Cells.Item(52,10).Value = 10
$ws1.Cells.Item(52,11).Value = 0.04
$ws1.Cells.Item(52,12).Value = 2800
$ws1.Cells.Item(52,13).Value = '55000-59500'
$ws1.Cells.Item(52,14).Value = '34000-40000'
$ws1.Cells.Item(52,15).Value = 1
$ws1.Cells.Item(52,18).Value = 'Labyrinth'

$ws1.Cells.Item(53,1).Value = 'Fabled Princess'
$ws1.Cells.Item(53,2).Value = 29000
$ws1.Cells.Item(53,3).Value = 28000
$ws1.Cells.Item(53,4).Value = 30600
$ws1.Cells.Item(53,5).Value = 27000
$ws1.Cells.Item(53,6).Value = 30000
$ws1.Cells.Item(53,7).Value = 850
$ws1.Cells.Item(53,8).Value = 525
$ws1.Cells.Item(53,9).Value = 'dex'
$ws1.Cells.Item(53,10).Value = 15
$ws1.Cells.Item(53,11).Value = 0.01
$ws1.Cells.Item(53,12).Value = 3000
$ws1.Cells.Item(53,13).Value = '60400-65700'
$ws1.Cells.Item(53,14).Value = '43000-50000'
$ws1.Cells.Item(53,15).Value = 1
$ws1.Cells.Item(53,16).Value = 'Weapon Smiths Book'
$ws1.Cells.Item(53,17).Value = 0.15
$ws1.Cells.Item(53,18).Value = 'Labyrinth'

$ws1.Cells.Item(54,1).Value = 'Demon Hunter'
$ws1.Cells.Item(54,2).Value = 29800
$ws1.Cells.Item(54,3).Value = 31000
$ws1.Cells.Item(54,4).Value = 30500
$ws1.Cells.Item(54,5).Value = 29800
$ws1.Cells.Item(54,6).Value = 2870
$ws1.Cells.Item(54,7).Value = 900
$ws1.Cells.Item(54,8).Value = 550
$ws1.Cells.Item(54,9).Value = 'dur'
$ws1.Cells.Item(54,10).Value = 15
$ws1.Cells.Item(54,11).Value = 0.04
$ws1.Cells.Item(54,12).Value = 3400
$ws1.Cells.Item(54,13).Value = '66700-69800'
$ws1.Cells.Item(54,14).Value = '66000-70000'
$ws1.Cells.Item(54,15).Value = 1
$ws1.Cells.Item(54,18).Value = 'Labyrinth'

$ws1.Cells.Item(55,1).Value = 'Cyclops Bat'
$ws1.Cells.Item(55,2).Value = 34000
$ws1.Cells.Item(55,3).Value = 35000
$ws1.Cells.Item(55,4).Value = 36000
$ws1.Cells.Item(55,5).Value = 32000
$ws1.Cells.Item(55,6).Value = 31000
$ws1.Cells.Item(55,7).Value = 930
$ws1.Cells.Item(55,8).Value = 625
$ws1.Cells.Item(55,9).Value = 'dex'
$ws1.Cells.Item(55,10).Value = 15
$ws1.Cells.Item(55,11).Value = 0.03
$ws1.Cells.Item(55,12).Value = 3500
$ws1.Cells.Item(55,13).Value = '70000-76000'
$ws1.Cells.Item(55,14).Value = '75000-83000'
$ws1.Cells.Item(55,15).Value = 1
$ws1.Cells.Item(55,18).Value = 'Labyrinth'

$ws1.Cells.Item(56,1).Value = 'Water Fiend'
$ws1.Cells.Item(56,2).Value = 36000
$ws1.Cells.Item(56,3).Value = 35000
$ws1.Cells.Item(56,4).Value = 38300
$ws1.Cells.Item(56,5).Value = 34000
$ws1.Cells.Item(56,6).Value = 34000
$ws1.Cells.Item(56,7).Value = 950
$ws1.Cells.Item(56,8).Value = 660
$ws1.Cells.Item(56,9).Value = 'dex'
$ws1.Cells.Item(56,10).Value = 15
$ws1.Cells.Item(56,11).Value = 0.03
$ws1.Cells.Item(56,12).Value = 3800
$ws1.Cells.Item(56,13).Value = '79500-86700'
$ws1.Cells.Item(56,14).Value = '84500-87000'
$ws1.Cells.Item(56,15).Value = 1
$ws1.Cells.Item(56,16).Value = 'Spell Weaving Book'
$ws1.Cells.Item(56,17).Value = 0.15
$ws1.Cells.Item(56,18).Value = 'Labyrinth'

$ws1.Cells.Item(57,1).Value = 'Labyrinth Hound'
$ws1.Cells.Item(57,2).Value = 38000
$ws1.Cells.Item(57,3).Value = 40000
$ws1.Cells.Item(57,4).Value = 39000
$ws1.Cells.Item(57,5).Value = 43000
$ws1.Cells.Item(57,6).Value = 43400
$ws1.Cells.Item(57,7).Value = 1000
$ws1.Cells.Item(57,8).Value = 750
$ws1.Cells.Item(57,9).Value = 'int'
$ws1.Cells.Item(57,10).Value = 20
$ws1.Cells.Item(57,11).Value = 0.04
$ws1.Cells.Item(57,12).Value = 4000
$ws1.Cells.Item(57,13).Value = '89000-94000'
$ws1.Cells.Item(57,14).Value = '90000-120000'
$ws1.Cells.Item(57,15).Value = 1
$ws1.Cells.Item(57,18).Value = 'Labyrinth'

# --- Sheet 'Monsters Skills': append new skill rows (84-113) ---
$ws2.Cells.Item(84,1).Value = 'Adult Red Dragon'
$ws2.Cells.Item(84,3).Value = 17
$ws2.Cells.Item(84,4).Value = 0
$ws2.Cells.Item(84,7).Value = 'Accuracy'

$ws2.Cells.Item(85,1).Value = 'Adult Red Dragon'
$ws2.Cells.Item(85,3).Value = 16
$ws2.Cells.Item(85,4).Value = 0
$ws2.Cells.Item(85,7).Value = 'Dodge'

$ws2.Cells.Item(86,1).Value = 'Litch Lord of the Labyrinth'
$ws2.Cells.Item(86,3).Value = 16
$ws2.Cells.Item(86,4).Value = 0
$ws2.Cells.Item(86,7).Value = 'Accuracy'

$ws2.Cells.Item(87,1).Value = 'Litch Lord of the Labyrinth'
$ws2.Cells.Item(87,3).Value = 18
$ws2.Cells.Item(87,4).Value = 0
$ws2.Cells.Item(87,7).Value = 'Dodge'

$ws2.Cells.Item(88,1).Value = 'Whisper of Fate'
$ws2.Cells.Item(88,3).Value = 19
$ws2.Cells.Item(88,4).Value = 0
$ws2.Cells.Item(88,7).Value = 'Accuracy'

$ws2.Cells.Item(89,1).Value = 'Whisper of Fate'
$ws2.Cells.Item(89,3).Value = 17
$ws2.Cells.Item(89,4).Value = 0
$ws2.Cells.Item(89,7).Value = 'Dodge'

$ws2.Cells.Item(90,1).Value = 'Demonic Angel'
$ws2.Cells.Item(90,3).Value = 20
$ws2.Cells.Item(90,4).Value = 0
$ws2.Cells.Item(90,7).Value = 'Accuracy'

$ws2.Cells.Item(91,1).Value = 'Demonic Angel'
$ws2.Cells.Item(91,3).Value = 17
$ws2.Cells.Item(91,4).Value = 0
$ws2.Cells.Item(91,7).Value = 'Dodge'

$ws2.Cells.Item(92,1).Value = 'Virgin Priestess Zombie'
$ws2.Cells.Item(92,3).Value = 21
$ws2.Cells.Item(92,4).Value = 0
$ws2.Cells.Item(92,7).Value = 'Accuracy'

$ws2.Cells.Item(93,1).Value = 'Virgin Priestess Zombie'
$ws2.Cells.Item(93,3).Value = 21
$ws2.Cells.Item(93,4).Value = 0
$ws2.Cells.Item(93,7).Value = 'Dodge'

$ws2.Cells.Item(94,1).Value = 'Hells Paladin'
$ws2.Cells.Item(94,3).Value = 20
$ws2.Cells.Item(94,4).Value = 0
$ws2.Cells.Item(94,7).Value = 'Accuracy'

$ws2.Cells.Item(95,1).Value = 'Hells Paladin'
$ws2.Cells.Item(95,3).Value = 20
$ws2.Cells.Item(95,4).Value = 0
$ws2.Cells.Item(95,7).Value = 'Dodge'

$ws2.Cells.Item(96,1).Value = 'Ancient Gold Dragon'
$ws2.Cells.Item(96,3).Value = 17
$ws2.Cells.Item(96,4).Value = 0
$ws2.Cells.Item(96,7).Value = 'Accuracy'

$ws2.Cells.Item(97,1).Value = 'Ancient Gold Dragon'
$ws2.Cells.Item(97,3).Value = 19
$ws2.Cells.Item(97,4).Value = 0
$ws2.Cells.Item(97,7).Value = 'Dodge'

$ws2.Cells.Item(98,1).Value = 'Broken Hearted Banshee'
$ws2.Cells.Item(98,3).Value = 24
$ws2.Cells.Item(98,4).Value = 0
$ws2.Cells.Item(98,7).Value = 'Accuracy'

$ws2.Cells.Item(99,1).Value = 'Broken Hearted Banshee'
$ws2.Cells.Item(99,3).Value = 23
$ws2.Cells.Item(99,4).Value = 0
$ws2.Cells.Item(99,7).Value = 'Dodge'

$ws2.Cells.Item(100,1).Value = 'Red Wizard of Krull'
$ws2.Cells.Item(100,3).Value = 21
$ws2.Cells.Item(100,4).Value = 0
$ws2.Cells.Item(100,7).Value = 'Accuracy'

$ws2.Cells.Item(101,1).Value = 'Red Wizard of Krull'
$ws2.Cells.Item(101,3).Value = 23
$ws2.Cells.Item(101,4).Value = 0
$ws2.Cells.Item(101,7).Value = 'Dodge'

$ws2.Cells.Item(102,1).Value = 'Devils Lover'
$ws2.Cells.Item(102,3).Value = 23
$ws2.Cells.Item(102,4).Value = 0
$ws2.Cells.Item(102,7).Value = 'Accuracy'

$ws2.Cells.Item(103,1).Value = 'Devils Lover'
$ws2.Cells.Item(103,3).Value = 22
$ws2.Cells.Item(103,4).Value = 0
$ws2.Cells.Item(103,7).Value = 'Dodge'

$ws2.Cells.Item(104,1).Value = 'Fabled Princess'
$ws2.Cells.Item(104,3).Value = 26
$ws2.Cells.Item(104,4).Value = 0
$ws2.Cells.Item(104,7).Value = 'Accuracy'

$ws2.Cells.Item(105,1).Value = 'Fabled Princess'
$ws2.Cells.Item(105,3).Value = 26
$ws2.Cells.Item(105,4).Value = 0
$ws2.Cells.Item(105,7).Value = 'Dodge'

$ws2.Cells.Item(106,1).Value = 'Demon Hunter'
$ws2.Cells.Item(106,3).Value = 24
$ws2.Cells.Item(106,4).Value = 0
$ws2.Cells.Item(106,7).Value = 'Accuracy'

$ws2.Cells.Item(107,1).Value = 'Demon Hunter'
$ws2.Cells.Item(107,3).Value = 23
$ws2.Cells.Item(107,4).Value = 0
$ws2.Cells.Item(107,7).Value = 'Dodge'

$ws2.Cells.Item(108,1).Value = 'Cyclops Bat'
$ws2.Cells.Item(108,3).Value = 20
$ws2.Cells.Item(108,4).Value = 0
$ws2.Cells.Item(108,7).Value = 'Accuracy'

$ws2.Cells.Item(109,1).Value = 'Cyclops Bat'
$ws2.Cells.Item(109,3).Value = 25
$ws2.Cells.Item(109,4).Value = 0
$ws2.Cells.Item(109,7).Value = 'Dodge'

$ws2.Cells.Item(110,1).Value = 'Water Fiend'
$ws2.Cells.Item(110,3).Value = 28
$ws2.Cells.Item(110,4).Value = 0
$ws2.Cells.Item(110,7).Value = 'Accuracy'

$ws2.Cells.Item(111,1).Value = 'Water Fiend'
$ws2.Cells.Item(111,3).Value = 25
$ws2.Cells.Item(111,4).Value = 0
$ws2.Cells.Item(111,7).Value = 'Dodge'

$ws2.Cells.Item(112,1).Value = 'Labyrinth Hound'
$ws2.Cells.Item(112,3).Value = 25
$ws2.Cells.Item(112,4).Value = 0
$ws2.Cells.Item(112,7).Value = 'Accuracy'

$ws2.Cells.Item(113,1).Value = 'Labyrinth Hound'
$ws2.Cells.Item(113,3).Value = 24
$ws2.Cells.Item(113,4).Value = 0
$ws2.Cells.Item(113,7).Value = 'Dodge'

# --- Column width adjustments ---
$offset = 0.8333333333333333
$ws1.Range("A1:A1").EntireColumn.ColumnWidth = 32 - $offset
$ws1.Range("B1:F1").EntireColumn.ColumnWidth = 6 - $offset
$ws1.Range("G1:G1").EntireColumn.ColumnWidth = 5 - $offset
$ws2.Range("A1:A1").EntireColumn.ColumnWidth = 32 - $offset
